# "small changes on laptop" - append one new trapping-visit record
# (year=2023, month=2, village=baiama, visit=9) as row 46 of Sheet1,
# and leave the selection where it lands after typing it (A47),
# mirroring the interactive edit captured in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(46, 1).Value = 2023
$ws.Cells.Item(46, 2).Value = 2
$ws.Cells.Item(46, 4).Value = "baiama"
$ws.Cells.Item(46, 5).Value = 9

# Move the active selection down to the next empty row, as Excel does
# after you finish typing a row and press Enter/Down.
$ws.Range("A47").Select() | Out-Null
